$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed data for rows 2-51 (columns B:E) from the updated crypto price feed
$data = @(
    @('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc','29.863.75','  +1.51%  '),
    @('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth','1.860.44','  +1.25%  '),
    @('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt','1.000','  +0.11%  '),
    @('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb','244.51','  +0.52%  '),
    @('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp','0.6555','  +4.88%  '),
    @('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc','1.001','  +0.07%  '),
    @('OKB','https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb','47.97','  +4.06%  '),
    @('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge','0.07575','  +2.45%  '),
    @('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada','0.2970','  +1.46%  '),
    @('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol','24.67','  +6.02%  '),
    @('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx','0.07653','  +0.22%  '),
    @('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth','1.874.65','  +2.04%  '),
    @('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot','5.079','  +1.41%  '),
    @('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic','0.6922','  +2.69%  '),
    @('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc','83.89','  +1.07%  '),
    @('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib','0.000009689','  +3.97%  '),
    @('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni','6.138','  +4.44%  '),
    @('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc','29.880.62','  +1.65%  '),
    @('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','2.121.01','  +1.62%  '),
    @('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch','236.62','  -0.65%  '),
    @('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax','12.71','  +1.76%  '),
    @('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai','1.001','  +0.06%  '),
    @('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link','7.791','  +6.06%  '),
    @('BinanceUSD','https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd','1.002','  +0.15%  '),
    @('Stellar','https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm','0.1441','  +1.92%  '),
    @('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr','158.72','  +0.09%  '),
    @('Cosmos','https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom','8.609','  +1.56%  '),
    @('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc','17.92','  +1.13%  '),
    @('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar','0.06060','  -1.37%  '),
    @('PancakeSwap','https://coinranking.com/coin/ncYFcP709+pancakeswap-cake','1.493','  -0.17%  '),
    @('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton','1.288','  +3.94%  '),
    @('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil','4.167','  +1.47%  '),
    @('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp','4.093','  +0.19%  '),
    @('LidoDAOToken','https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo','1.876','  +0.57%  '),
    @('ARBITRUM','https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb','1.182','  +3.67%  '),
    @('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx','0.7357','  +1.45%  '),
    @('HuobiToken','https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht','2.611','  +0.13%  '),
    @('MXToken','https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx','2.817','  -2.39%  '),
    @('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet','0.01794','  +1.83%  '),
    @('Maker','https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr','1.208.50','  -0.67%  '),
    @('FraxShare','https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs','6.326','  +0.50%  '),
    @('TrustWalletToken','https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt','0.9157','  +0.25%  '),
    @('RocketPoolETH','https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth','2.036.22','  +1.75%  '),
    @('PaxDollar','https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp','1.001','  -0.03%  '),
    @('Aave','https://coinranking.com/coin/ixgUfzmLR+aave-aave','67.46','  +3.24%  '),
    @('Quant','https://coinranking.com/coin/bauj_21eYVwso+quant-qnt','101.64','  +0.03%  '),
    @('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt','7.502','  +12.74%  '),
    @('BabyDogeCoin','https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge','0.00000000121','  +2.30%  '),
    @('TheSandbox','https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand','0.4075','  +0.61%  '),
    @('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens','9.225','  +0.03%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($rowNum, 2).Value = $vals[0]
    $ws.Cells.Item($rowNum, 3).Value = $vals[1]

    $priceCell = $ws.Cells.Item($rowNum, 4)
    $priceVal = $vals[2]
    $looksNumeric = $priceVal -match "^[0-9]*\.?[0-9]+$"
    if ($looksNumeric) {
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceVal
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $priceVal
    }

    $ws.Cells.Item($rowNum, 5).Value = $vals[3]
}
